$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-08-17T19:07:03"
$ws.Range("W4").Value = 90
$ws.Range("X4").Value = 52.63
$ws.Range("Y4").Value = 34.33
$ws.Range("Z4").Value = 34.16
$ws.Range("W6").Value = -2.43
$ws.Range("X6").Value = -1
$ws.Range("Y6").Value = -0.31
$ws.Range("Z6").Value = -0.58
$ws.Range("W9").Value = 93.56
$ws.Range("X9").Value = 75
$ws.Range("Y9").Value = 35.49
$ws.Range("Z9").Value = 35.56
$ws.Range("W11").Value = 0.93
$ws.Range("X11").Value = 0.98
$ws.Range("Y11").Value = 0.85
$ws.Range("Z11").Value = 0.82
$ws.Range("W12").Value = 0.19
$ws.Range("X12").Value = 20.39
$ws.Range("W14").Value = 93.65000000000001
$ws.Range("X14").Value = 75
$ws.Range("Y14").Value = 35.53
$ws.Range("Z14").Value = 35.56
$ws.Range("W16").Value = 1.03
$ws.Range("X16").Value = 0.98
$ws.Range("Y16").Value = 0.89
$ws.Range("Z16").Value = 0.82
$ws.Range("W17").Value = 0.19
$ws.Range("X17").Value = 20.39
$ws.Range("W19").Value = 89.81999999999999
$ws.Range("X19").Value = 52.48
$ws.Range("Y19").Value = 34.23
$ws.Range("Z19").Value = 34.06
$ws.Range("W21").Value = -2.6
$ws.Range("X21").Value = -1.15
$ws.Range("Y21").Value = -0.41
$ws.Range("Z21").Value = -0.68
$ws.Range("W24").Value = 89.81999999999999
$ws.Range("X24").Value = 52.48
$ws.Range("Y24").Value = 34.23
$ws.Range("Z24").Value = 34.06
$ws.Range("W26").Value = -2.6
$ws.Range("X26").Value = -1.15
$ws.Range("Y26").Value = -0.41
$ws.Range("Z26").Value = -0.68
$ws.Range("W29").Value = 89.3
$ws.Range("X29").Value = 52.12
$ws.Range("Y29").Value = 33.99
$ws.Range("Z29").Value = 33.8
$ws.Range("W31").Value = -3.13
$ws.Range("X31").Value = -1.51
$ws.Range("Y31").Value = -0.65
$ws.Range("Z31").Value = -0.95
$ws.Range("W34").Value = 96.18000000000001
$ws.Range("X34").Value = 56.82
$ws.Range("Y34").Value = 36.5
$ws.Range("Z34").Value = 36.69
$ws.Range("W36").Value = 3.75
$ws.Range("X36").Value = 2.7
$ws.Range("Y36").Value = 1.86
$ws.Range("Z36").Value = 1.94
$ws.Range("X37").Value = 0.49
$ws.Range("W39").Value = 90
$ws.Range("X39").Value = 52.63
$ws.Range("Y39").Value = 34.33
$ws.Range("Z39").Value = 34.16
$ws.Range("W41").Value = -2.43
$ws.Range("X41").Value = -1
$ws.Range("Y41").Value = -0.31
$ws.Range("Z41").Value = -0.58
$ws.Range("W44").Value = 93
$ws.Range("X44").Value = 54.29
$ws.Range("Y44").Value = 35.13
$ws.Range("Z44").Value = 35.07
$ws.Range("W46").Value = 0.57
$ws.Range("X46").Value = 0.66
$ws.Range("Y46").Value = 0.49
$ws.Range("Z46").Value = 0.33
$ws.Range("W49").Value = 95.19
$ws.Range("X49").Value = 55.4
$ws.Range("Y49").Value = 35.49
$ws.Range("Z49").Value = 35.75
$ws.Range("W51").Value = 2.76
$ws.Range("X51").Value = 1.77
$ws.Range("Y51").Value = 0.85
$ws.Range("Z51").Value = 1
$ws.Range("W54").Value = 89.47
$ws.Range("X54").Value = 52.52
$ws.Range("Y54").Value = 34.37
$ws.Range("Z54").Value = 34.61
$ws.Range("W56").Value = -2.78
$ws.Range("X56").Value = -0.79
$ws.Range("Y56").Value = -0.27
$ws.Range("Z56").Value = -0.14
$ws.Range("W57").Value = -0.18
$ws.Range("X57").Value = -0.32
$ws.Range("W59").Value = 96.68000000000001
$ws.Range("X59").Value = 56.16
$ws.Range("Y59").Value = 36.12
$ws.Range("Z59").Value = 36.08
$ws.Range("W61").Value = 4.25
$ws.Range("X61").Value = 2.53
$ws.Range("Y61").Value = 1.48
$ws.Range("Z61").Value = 1.33
$ws.Range("W64").Value = 97.91
$ws.Range("X64").Value = 56.81
$ws.Range("Y64").Value = 36.5
$ws.Range("Z64").Value = 36.46
$ws.Range("W66").Value = 5.48
$ws.Range("X66").Value = 3.18
$ws.Range("Y66").Value = 1.86
$ws.Range("Z66").Value = 1.71
$ws.Range("W69").Value = 98.75
$ws.Range("X69").Value = 57.36
$ws.Range("Y69").Value = 36.85
$ws.Range("Z69").Value = 36.92
$ws.Range("W71").Value = 6.32
$ws.Range("X71").Value = 3.73
$ws.Range("Y71").Value = 2.21
$ws.Range("Z71").Value = 2.18
$ws.Range("W74").Value = 95.98
$ws.Range("X74").Value = 55.69
$ws.Range("Y74").Value = 35.79
$ws.Range("Z74").Value = 35.78
$ws.Range("W76").Value = 3.55
$ws.Range("X76").Value = 2.06
$ws.Range("Y76").Value = 1.15
$ws.Range("Z76").Value = 1.04
$ws.Range("W79").Value = 96.63
$ws.Range("X79").Value = 56.07
$ws.Range("Y79").Value = 36.06
$ws.Range("Z79").Value = 36.04
$ws.Range("W81").Value = 4.2
$ws.Range("X81").Value = 2.44
$ws.Range("Y81").Value = 1.42
$ws.Range("Z81").Value = 1.29
$ws.Range("W84").Value = 90.26000000000001
$ws.Range("X84").Value = 53
$ws.Range("Y84").Value = 34.26
$ws.Range("Z84").Value = 34.43
$ws.Range("W86").Value = -2.17
$ws.Range("X86").Value = -0.64
$ws.Range("Y86").Value = -0.38
$ws.Range("Z86").Value = -0.31
$ws.Range("W89").Value = 89.22
$ws.Range("X89").Value = 52.12
$ws.Range("Y89").Value = 33.99
$ws.Range("Z89").Value = 33.8
$ws.Range("W91").Value = -3.21
$ws.Range("X91").Value = -1.51
$ws.Range("Y91").Value = -0.65
$ws.Range("Z91").Value = -0.95
